$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Data")

# New report date: 2026-01-29 (Excel serial 46051), appended after existing row 397 (2026-01-28)
$dateSerial = 46051

# Pre-format the new date column cells to match the existing date column style/format
$ws.Range("A398:A419").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(398, 1).Value = $dateSerial
$ws.Cells.Item(398, 2).Value = 'ASAHI DEPOSITORY LLC Registered'
$ws.Cells.Item(398, 3).Value = 0
$ws.Cells.Item(398, 4).Value = 0
$ws.Cells.Item(398, 5).Value = 0
$ws.Cells.Item(398, 6).Value = 0
$ws.Cells.Item(398, 7).Value = 0
$ws.Cells.Item(398, 8).Value = 0

$ws.Cells.Item(399, 1).Value = $dateSerial
$ws.Cells.Item(399, 2).Value = 'ASAHI DEPOSITORY LLC Eligible'
$ws.Cells.Item(399, 3).Value = 0
$ws.Cells.Item(399, 4).Value = 0
$ws.Cells.Item(399, 5).Value = 0
$ws.Cells.Item(399, 6).Value = 0
$ws.Cells.Item(399, 7).Value = 0
$ws.Cells.Item(399, 8).Value = 0

$ws.Cells.Item(400, 1).Value = $dateSerial
$ws.Cells.Item(400, 2).Value = 'BRINK''S, INC. Registered'
$ws.Cells.Item(400, 3).Value = 87949.747
$ws.Cells.Item(400, 4).Value = 0
$ws.Cells.Item(400, 5).Value = 0
$ws.Cells.Item(400, 6).Value = 0
$ws.Cells.Item(400, 7).Value = 0
$ws.Cells.Item(400, 8).Value = 87949.747

$ws.Cells.Item(401, 1).Value = $dateSerial
$ws.Cells.Item(401, 2).Value = 'BRINK''S, INC. Eligible'
$ws.Cells.Item(401, 3).Value = 30578.352
$ws.Cells.Item(401, 4).Value = 0
$ws.Cells.Item(401, 5).Value = 0
$ws.Cells.Item(401, 6).Value = 0
$ws.Cells.Item(401, 7).Value = 0
$ws.Cells.Item(401, 8).Value = 30578.352

$ws.Cells.Item(402, 1).Value = $dateSerial
$ws.Cells.Item(402, 2).Value = 'CNT DEPOSITORY, INC. Registered'
$ws.Cells.Item(402, 3).Value = 1246.06
$ws.Cells.Item(402, 4).Value = 0
$ws.Cells.Item(402, 5).Value = 0
$ws.Cells.Item(402, 6).Value = 0
$ws.Cells.Item(402, 7).Value = 0
$ws.Cells.Item(402, 8).Value = 1246.06

$ws.Cells.Item(403, 1).Value = $dateSerial
$ws.Cells.Item(403, 2).Value = 'CNT DEPOSITORY, INC. Eligible'
$ws.Cells.Item(403, 3).Value = 0
$ws.Cells.Item(403, 4).Value = 0
$ws.Cells.Item(403, 5).Value = 0
$ws.Cells.Item(403, 6).Value = 0
$ws.Cells.Item(403, 7).Value = 0
$ws.Cells.Item(403, 8).Value = 0

$ws.Cells.Item(404, 1).Value = $dateSerial
$ws.Cells.Item(404, 2).Value = 'DELAWARE DEPOSITORY Registered'
$ws.Cells.Item(404, 3).Value = 1633.941
$ws.Cells.Item(404, 4).Value = 0
$ws.Cells.Item(404, 5).Value = 0
$ws.Cells.Item(404, 6).Value = 0
$ws.Cells.Item(404, 7).Value = 0
$ws.Cells.Item(404, 8).Value = 1633.941

$ws.Cells.Item(405, 1).Value = $dateSerial
$ws.Cells.Item(405, 2).Value = 'DELAWARE DEPOSITORY Eligible'
$ws.Cells.Item(405, 3).Value = 18459.584
$ws.Cells.Item(405, 4).Value = 0
$ws.Cells.Item(405, 5).Value = 0
$ws.Cells.Item(405, 6).Value = 0
$ws.Cells.Item(405, 7).Value = 0
$ws.Cells.Item(405, 8).Value = 18459.584

$ws.Cells.Item(406, 1).Value = $dateSerial
$ws.Cells.Item(406, 2).Value = 'HSBC BANK, USA Registered'
$ws.Cells.Item(406, 3).Value = 1394.758
$ws.Cells.Item(406, 4).Value = 0
$ws.Cells.Item(406, 5).Value = 0
$ws.Cells.Item(406, 6).Value = 0
$ws.Cells.Item(406, 7).Value = 0
$ws.Cells.Item(406, 8).Value = 1394.758

$ws.Cells.Item(407, 1).Value = $dateSerial
$ws.Cells.Item(407, 2).Value = 'HSBC BANK, USA Eligible'
$ws.Cells.Item(407, 3).Value = 9281.978999999999
$ws.Cells.Item(407, 4).Value = 0
$ws.Cells.Item(407, 5).Value = 0
$ws.Cells.Item(407, 6).Value = 0
$ws.Cells.Item(407, 7).Value = 0
$ws.Cells.Item(407, 8).Value = 9281.978999999999

$ws.Cells.Item(408, 1).Value = $dateSerial
$ws.Cells.Item(408, 2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered'
$ws.Cells.Item(408, 3).Value = 2395.448
$ws.Cells.Item(408, 4).Value = 0
$ws.Cells.Item(408, 5).Value = 0
$ws.Cells.Item(408, 6).Value = 0
$ws.Cells.Item(408, 7).Value = 0
$ws.Cells.Item(408, 8).Value = 2395.448

$ws.Cells.Item(409, 1).Value = $dateSerial
$ws.Cells.Item(409, 2).Value = 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible'
$ws.Cells.Item(409, 3).Value = 0
$ws.Cells.Item(409, 4).Value = 0
$ws.Cells.Item(409, 5).Value = 0
$ws.Cells.Item(409, 6).Value = 0
$ws.Cells.Item(409, 7).Value = 0
$ws.Cells.Item(409, 8).Value = 0

$ws.Cells.Item(410, 1).Value = $dateSerial
$ws.Cells.Item(410, 2).Value = 'JP MORGAN CHASE BANK NA Registered'
$ws.Cells.Item(410, 3).Value = 114985.579
$ws.Cells.Item(410, 4).Value = 0
$ws.Cells.Item(410, 5).Value = 0
$ws.Cells.Item(410, 6).Value = 0
$ws.Cells.Item(410, 7).Value = 0
$ws.Cells.Item(410, 8).Value = 114985.579

$ws.Cells.Item(411, 1).Value = $dateSerial
$ws.Cells.Item(411, 2).Value = 'JP MORGAN CHASE BANK NA Eligible'
$ws.Cells.Item(411, 3).Value = 125407.673
$ws.Cells.Item(411, 4).Value = 0
$ws.Cells.Item(411, 5).Value = 0
$ws.Cells.Item(411, 6).Value = 0
$ws.Cells.Item(411, 7).Value = 0
$ws.Cells.Item(411, 8).Value = 125407.673

$ws.Cells.Item(412, 1).Value = $dateSerial
$ws.Cells.Item(412, 2).Value = 'LOOMIS INTERNATIONAL (US) LLC Registered'
$ws.Cells.Item(412, 3).Value = 63745.991
$ws.Cells.Item(412, 4).Value = 0
$ws.Cells.Item(412, 5).Value = 0
$ws.Cells.Item(412, 6).Value = 0
$ws.Cells.Item(412, 7).Value = 0
$ws.Cells.Item(412, 8).Value = 63745.991

$ws.Cells.Item(413, 1).Value = $dateSerial
$ws.Cells.Item(413, 2).Value = 'LOOMIS INTERNATIONAL (US) LLC Eligible'
$ws.Cells.Item(413, 3).Value = 132077.206
$ws.Cells.Item(413, 4).Value = 0
$ws.Cells.Item(413, 5).Value = 0
$ws.Cells.Item(413, 6).Value = 0
$ws.Cells.Item(413, 7).Value = 0
$ws.Cells.Item(413, 8).Value = 132077.206

$ws.Cells.Item(414, 1).Value = $dateSerial
$ws.Cells.Item(414, 2).Value = 'MALCA-AMIT USA, LLC Registered'
$ws.Cells.Item(414, 3).Value = 395.145
$ws.Cells.Item(414, 4).Value = 0
$ws.Cells.Item(414, 5).Value = 0
$ws.Cells.Item(414, 6).Value = 0
$ws.Cells.Item(414, 7).Value = 0
$ws.Cells.Item(414, 8).Value = 395.145

$ws.Cells.Item(415, 1).Value = $dateSerial
$ws.Cells.Item(415, 2).Value = 'MALCA-AMIT USA, LLC Eligible'
$ws.Cells.Item(415, 3).Value = 0
$ws.Cells.Item(415, 4).Value = 0
$ws.Cells.Item(415, 5).Value = 0
$ws.Cells.Item(415, 6).Value = 0
$ws.Cells.Item(415, 7).Value = 0
$ws.Cells.Item(415, 8).Value = 0

$ws.Cells.Item(416, 1).Value = $dateSerial
$ws.Cells.Item(416, 2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Registered'
$ws.Cells.Item(416, 3).Value = 50220.42
$ws.Cells.Item(416, 4).Value = 0
$ws.Cells.Item(416, 5).Value = 0
$ws.Cells.Item(416, 6).Value = 0
$ws.Cells.Item(416, 7).Value = 0
$ws.Cells.Item(416, 8).Value = 50220.42

$ws.Cells.Item(417, 1).Value = $dateSerial
$ws.Cells.Item(417, 2).Value = 'MANFRA, TORDELLA & BROOKES, LLC Eligible'
$ws.Cells.Item(417, 3).Value = 1271.373
$ws.Cells.Item(417, 4).Value = 0
$ws.Cells.Item(417, 5).Value = 0
$ws.Cells.Item(417, 6).Value = 0
$ws.Cells.Item(417, 7).Value = 0
$ws.Cells.Item(417, 8).Value = 1271.373

$ws.Cells.Item(418, 1).Value = $dateSerial
$ws.Cells.Item(418, 2).Value = 'STONEX PRECIOUS METALS LLC Registered'
$ws.Cells.Item(418, 3).Value = 14122.765
$ws.Cells.Item(418, 4).Value = 0
$ws.Cells.Item(418, 5).Value = 0
$ws.Cells.Item(418, 6).Value = 0
$ws.Cells.Item(418, 7).Value = 0
$ws.Cells.Item(418, 8).Value = 14122.765

$ws.Cells.Item(419, 1).Value = $dateSerial
$ws.Cells.Item(419, 2).Value = 'STONEX PRECIOUS METALS LLC Eligible'
$ws.Cells.Item(419, 3).Value = 16.075
$ws.Cells.Item(419, 4).Value = 0
$ws.Cells.Item(419, 5).Value = 0
$ws.Cells.Item(419, 6).Value = 0
$ws.Cells.Item(419, 7).Value = 0
$ws.Cells.Item(419, 8).Value = 16.075

"Appended rows 398-419 (date 2026-01-29) to Daily_Data"